$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NATMI re-run following Dr Hou's advice: "ECs" is added as a sending/target
# cluster alongside the existing "FAPs" and "sCs" clusters. The Sending
# cluster x Target cluster combinations are now a full 3x3 cross join
# (previously 2 senders x 3 targets = 6 rows), so the data block grows from
# rows 2-7 to rows 2-10 and the used range grows from A1:T7 to A1:T10.
#
# Columns (A-T), matching the existing header row (row 1) order:
#   A Sending cluster            K Receptor-expressing cells
#   B Ligand symbol              L Receptor detection rate
#   C Receptor symbol            M Receptor average expression value
#   D Target cluster             N Receptor total expression value
#   E Ligand-expressing cells    O Receptor derived specificity (avg)
#   F Ligand detection rate      P Receptor derived specificity (total)
#   G Ligand average expr value  Q Edge average expression weight
#   H Ligand total expr value    R Edge total expression weight
#   I Ligand derived spec (avg)  S Edge average expression derived specificity
#   J Ligand derived spec (tot)  T Edge total expression derived specificity
$data = @(
    @("ECs",  "Fgf18", "Fgfr1", "ECs",  1, 0.3333333333333333, 0.4765323333333333, 1.429597,   0.04733670697480491, 0.04733670697480491, 3, 1, 5.900730666666667,  17.702192,   0.03970749001357476, 0.03970749001357476, 2.811888952958222,  25.307000576624,   0.001879621819477581, 0.001879621819477581),
    @("ECs",  "Fgf18", "Fgfr1", "FAPs", 1, 0.3333333333333333, 0.4765323333333333, 1.429597,   0.04733670697480491, 0.04733670697480491, 3, 1, 111.5917106666667,  334.775132,  0.7509284844884279,  0.7509284844884279,  53.17705826464488,  478.593524381804,  0.03554648162926304,  0.03554648162926304),
    @("ECs",  "Fgf18", "Fgfr1", "sCs",  1, 0.3333333333333333, 0.4765323333333333, 1.429597,   0.04733670697480491, 0.04733670697480491, 3, 1, 31.11253633333333,  93.337609,   0.2093640254979974,  0.2093640254979974,  14.82612953484145,  133.435165813573,  0.009910603526064288, 0.009910603526064286),
    @("FAPs", "Fgf18", "Fgfr1", "ECs",  3, 1,                   8.770835333333332,  26.312506,  0.8712576945074703,  0.8712576945074703,  3, 1, 5.900730666666667,  17.702192,   0.03970749001357476, 0.03970749001357476, 51.75433702368355,  465.789033213152,  0.03459545620390555,  0.03459545620390554),
    @("FAPs", "Fgf18", "Fgfr1", "FAPs", 3, 1,                   8.770835333333332,  26.312506,  0.8712576945074703,  0.8712576945074703,  3, 1, 111.5917106666667,  334.775132,  0.7509284844884279,  0.7509284844884279,  978.75251882231,    8808.77266940079,  0.6542522201353763,   0.6542522201353763),
    @("FAPs", "Fgf18", "Fgfr1", "sCs",  3, 1,                   8.770835333333332,  26.312506,  0.8712576945074703,  0.8712576945074703,  3, 1, 31.11253633333333,  93.337609,   0.2093640254979974,  0.2093640254979974,  272.8829329820171,  2455.946396838154, 0.1824100181681885,   0.1824100181681884),
    @("sCs",  "Fgf18", "Fgfr1", "ECs",  3, 1,                   0.8194993333333332, 2.458498,   0.08140559851772487, 0.08140559851772486, 3, 1, 5.900730666666667,  17.702192,   0.03970749001357476, 0.03970749001357476, 4.835644847512889,  43.520803627616,   0.003232411990191637, 0.003232411990191636),
    @("sCs",  "Fgf18", "Fgfr1", "FAPs", 3, 1,                   0.8194993333333332, 2.458498,   0.08140559851772487, 0.08140559851772486, 3, 1, 111.5917106666667,  334.775132,  0.7509284844884279,  0.7509284844884279,  91.44933249685954,  823.0439924717358, 0.06112978272378855,  0.06112978272378854),
    @("sCs",  "Fgf18", "Fgfr1", "sCs",  3, 1,                   0.8194993333333332, 2.458498,   0.08140559851772487, 0.08140559851772486, 3, 1, 31.11253633333333,  93.337609,   0.2093640254979974,  0.2093640254979974,  25.49670278347578,  229.470325051282,  0.01704340380374469,  0.01704340380374469)
)

$rowIndex = 2
foreach ($row in $data) {
    $colIndex = 1
    foreach ($val in $row) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $val
        $colIndex = $colIndex + 1
    }
    $rowIndex = $rowIndex + 1
}
